$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 6 (keep header row 1 and first data row 2)
$ws.Rows("3:6").Delete()

# Update the first data row with the new value/date
$ws.Range("A2").Value = "2015-11-18 18:10"
$ws.Range("B2").Value = 3.87
